# Update imputed values in columns A and C (result_data_RandomForest)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.69999999999999
$ws.Range("A7").Value = -19.62709999999999
$ws.Range("C7").Value = -12.62740000000001
$ws.Range("C15").Value = -14.54439999999999
$ws.Range("A16").Value = -21.63139999999999
$ws.Range("C21").Value = -12.2554
$ws.Range("C22").Value = -12.62660000000001
$ws.Range("C23").Value = -12.74730000000001
$ws.Range("A28").Value = -22.16239999999999
$ws.Range("A29").Value = -21.23929999999999
$ws.Range("A32").Value = -21.1808
$ws.Range("C34").Value = -11.26510000000002
$ws.Range("A40").Value = -20.15239999999999
$ws.Range("C43").Value = -12.65759999999999
$ws.Range("C45").Value = -13.35819999999998
$ws.Range("C50").Value = -13.91879999999998
$ws.Range("C51").Value = -11.9811
$ws.Range("A52").Value = -22.2623
$ws.Range("A57").Value = -22.2625
$ws.Range("A66").Value = -21.9591
$ws.Range("C66").Value = -12.2387
$ws.Range("C67").Value = -10.9857
$ws.Range("C79").Value = -11.3472
$ws.Range("C84").Value = -13.26929999999999
$ws.Range("C92").Value = -11.4224
$ws.Range("C97").Value = -12.3951
$ws.Range("A100").Value = -21.8824
